# Updated cryptos list on Sat Dec  9 15:35:07 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) / Volume(1h) (E) figures for each coin row, and
# re-sorts three adjacent row pairs whose relative ranking flipped
# (Stellar/Kaspa, Cronos/InjectiveProtocol, Aave/ARBITRUM) by swapping the
# Coin/Link/Price/Volume cell contents between the two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text updates (values that Excel would never mistake for a number,
# e.g. two-dot thousands separators, percentages, or coin names/links).
$textUpdates = [ordered]@{
    "D2"  = "44.006.35";  "E2"  = "  +0.12%  "
    "D3"  = "2.354.59";   "E3"  = "  -0.33%  "
    "E4"  = "  +0.19%  "
    "E5"  = "  +3.45%  "
    "E6"  = "  +1.81%  "
    "E7"  = "  +4.56%  "
    "E9"  = "  +14.07%  "
    "E10" = "  +2.20%  "
    "E11" = "  +0.22%  "
    "E12" = "  +17.21%  "
    "E13" = "  +11.58%  "
    "E14" = "  +1.27%  "
    "D15" = "2.702.46";   "E15" = "  -0.47%  "
    "E16" = "  -1.48%  "
    "E17" = "  +3.05%  "
    "D18" = "2.340.66";   "E18" = "  -1.06%  "
    "D19" = "43.859.91";  "E19" = "  -0.10%  "
    "E20" = "  +0.96%  "
    "E21" = "  +4.31%  "
    "E22" = "  +1.29%  "
    "E23" = "  +2.55%  "
    "E24" = "  +0.05%  "
    "E25" = "  -2.19%  "
    "E27" = "  +16.15%  "
    "E28" = "  +3.80%  "
    "E29" = "  +1.63%  "
    "E30" = "  -0.85%  "
    "E31" = "  +1.00%  "
    "B32" = "Kaspa";  "C32" = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas";                     "E32" = "  -4.22%  "
    "B33" = "Stellar"; "C33" = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm";               "E33" = "  +4.34%  "
    "E34" = "  +2.25%  "
    "E35" = "  +6.04%  "
    "E36" = "  +4.16%  "
    "E37" = "  -0.88%  "
    "E38" = "  -1.93%  "
    "E39" = "  -0.82%  "
    "E40" = "  +5.44%  "
    "E41" = "  +17.80%  "
    "B42" = "InjectiveProtocol"; "C42" = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"; "E42" = "  -0.58%  "
    "B43" = "Cronos";            "C43" = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro";          "E43" = "  +10.84%  "
    "E44" = "  +1.66%  "
    "E45" = "  -0.04%  "
    "E46" = "  +7.64%  "
    "E47" = "  +9.03%  "
    "E48" = "  +2.97%  "
    "B49" = "ARBITRUM"; "C49" = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"; "E49" = "  +0.24%  "
    "B50" = "Aave";     "C50" = "https://coinranking.com/coin/ixgUfzmLR+aave-aave";    "E50" = "  +2.27%  "
    "E51" = "  +7.15%  "
}

foreach ($ref in $textUpdates.Keys) {
    $ws.Range($ref).Value = $textUpdates[$ref]
}

# Price updates whose new text would otherwise be auto-parsed by Excel as a
# number (e.g. "239.52", "2.50"), silently dropping the original text
# formatting. Force the cell to Text format first so the literal string is
# preserved exactly as in the source feed, then drop the forced format back
# to the sheet's default style (only the *value* should differ from before).
$numericLookingUpdates = [ordered]@{
    "D6"  = "239.52"
    "D7"  = "75.90"
    "D9"  = "0.615"
    "D11" = "57.17"
    "D12" = "33.12"
    "D13" = "7.41"
    "D14" = "0.107"
    "D16" = "16.59"
    "D21" = "6.61"
    "D22" = "77.28"
    "D23" = "258.04"
    "D26" = "2.50"
    "D27" = "1.81"
    "D28" = "10.76"
    "D29" = "22.89"
    "D31" = "174.64"
    "D32" = "0.127"
    "D33" = "0.137"
    "D39" = "6.36"
    "D41" = "0.213"
    "D42" = "19.22"
    "D43" = "0.108"
    "D44" = "9.09"
    "D49" = "1.18"
    "D50" = "100.32"
    "D51" = "55.84"
}

foreach ($ref in $numericLookingUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $numericLookingUpdates[$ref]
    $cell.Style = "Normal"
}
